# Applies the text edits from the commit:
#   Slide 3 ("Technologies et outils" / "Content Placeholder 2"):
#     - "Langage : C#"   -> "Langage : C# Langage de programmation"
#     - "Framework : WPF" -> "Framework : WPF Pour l'interface visuel"
#     (the other bullets on this slide keep the same text; only a
#      cosmetic/spell-check "dirty" run flag changed for them upstream,
#      which carries no visible/textual effect)
#   Slide 7 ("Fonctionnalités futures" / "Content Placeholder 2"):
#     - merges the two runs "Création de " + "comptes par l'admin" of the
#       last bullet into a single run "Création de comptes par l'admin"

$p = $ppt.ActivePresentation

# --- Slide 3: shape 3 "Content Placeholder 2" ---
$s3 = $p.Slides.Item(3)
$sh3 = $s3.Shapes.Item(3)
$tr3 = $sh3.TextFrame.TextRange

$tr3.Paragraphs(1).Runs(1).Text = "Langage : C# Langage de programmation"
$tr3.Paragraphs(2).Runs(1).Text = "Framework : WPF Pour l’interface visuel"

# --- Slide 7: shape 2 "Content Placeholder 2" ---
$s7 = $p.Slides.Item(7)
$sh7 = $s7.Shapes.Item(2)
$tr7 = $sh7.TextFrame.TextRange

# Last paragraph ("Création de " + "comptes par l'admin") -> single run.
$para4 = $tr7.Paragraphs(4)
$run1 = $para4.Runs(1)
$run2 = $para4.Runs(2)
$run1.Text = "Création de comptes par l’admin"
$run2.Text = ""
